$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (device master test data) - rows 157..161
$ids     = @(3000176, 3000177, 3000178, 3000179, 3000180)
$names   = @("Finger Print Scanner 32", "IRIS Scanner 32", "Web Camera 32", "Document Scanner 32", "Printer 32")
$macs    = @("80-75-40-E8-CA-24", "0E-1A-14-4A-6D-3A", "65-13-7F-0F-F7-53", "73-C4-DE-8E-C9-8D", "EC-74-AB-E0-0F-38")
$serials = @("BS563Q2230824", "BS563Q2230825", "BS563Q2230826", "BS563Q2230827", "BS563Q2230828")
$dspecs  = @(165, 327, 736, 801, 920)

$startRow = 157

# Write column by column so new shared-string entries are appended in the
# same order as the target workbook (all names, then all macs, then all
# serials, then the repeated lang/active/cr_by/cr_dtimes values).
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $ids[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $names[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $macs[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $serials[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 6).Value = $dspecs[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 7).Value = "eng"
}
for ($i = 0; $i -lt 5; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 8).Value = $true
    $ws.Cells.Item($r, 8).HorizontalAlignment = -4131
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 9).Value = "superadmin"
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 10).Value = "now()"
}

# Update the view: scrolled position and selection, matching the saved
# workbook view state (entire columns K:XFD selected, scrolled to row 113).
$ws.Select()
$ws.Range("K1:XFD1048576").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 113
$win.ScrollColumn = 1
